# Generated PowerShell-style Excel COM-interop edit script
# Applies the weekly CompStat crime-data refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: "Volume 31   Number  15" -> "...Number  16" ---
$ws.Range("A8").Value = "Volume 31   Number  16"

# --- Header: report week date range ---
$ws.Range("C9").Value = "Report Covering the Week  4/15/2024  Through  4/21/2024"

# --- Stable reference cells (outside the edited row range) used to copy the
#     cell format+value for cells whose underlying value type flips between a
#     text placeholder ("0" / "***.*") and a real number. Copying cell-to-cell
#     transfers both the style index and the value/shared-string in one shot. ---
# C14 -> General style (s=14) holding shared string "0"
# E14 -> General style (s=14) holding shared string "***.*"
# J14 -> Number style (s=15) holding a plain number

# Row 15
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("L15").Value = 66.666666666666

# Row 16
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 400
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 180
$ws.Range("I16").Value = 51
$ws.Range("J16").Value = 33
$ws.Range("K16").Value = 54.545454545454
$ws.Range("L16").Value = 15.909090909090
$ws.Range("M16").Value = -15
$ws.Range("N16").Value = -79.012345679012

# Row 17
$ws.Range("J14").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = -50
$ws.Range("I17").Value = 38
$ws.Range("J17").Value = 45
$ws.Range("K17").Value = -15.555555555555
$ws.Range("L17").Value = -13.636363636363
$ws.Range("M17").Value = 15.151515151515
$ws.Range("N17").Value = -67.521367521367

# Row 18
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -40
$ws.Range("I18").Value = 38
$ws.Range("K18").Value = -11.627906976744
$ws.Range("L18").Value = -48.648648648648
$ws.Range("M18").Value = -2.564102564102
$ws.Range("N18").Value = -88.985507246376

# Row 19
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 14.285714285714
$ws.Range("I19").Value = 126
$ws.Range("J19").Value = 155
$ws.Range("K19").Value = -18.709677419354
$ws.Range("L19").Value = -20.253164556962
$ws.Range("M19").Value = -16
$ws.Range("N19").Value = -59.090909090909

# Row 20
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 3
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -72.727272727272
$ws.Range("I20").Value = 15
$ws.Range("J20").Value = 34
$ws.Range("K20").Value = -55.882352941176
$ws.Range("L20").Value = -16.666666666666
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -94.773519163763

# Row 21
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 13.333333333333
$ws.Range("F21").Value = 66
$ws.Range("G21").Value = 74
$ws.Range("H21").Value = -10.810810810810
$ws.Range("I21").Value = 273
$ws.Range("J21").Value = 317
$ws.Range("K21").Value = -13.880126182965
$ws.Range("L21").Value = -19.941348973607
$ws.Range("M21").Value = -7.770270270270
$ws.Range("N21").Value = -79.271070615034

# Row 22
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("M22").Value = -9.090909090909

# Row 23
$ws.Range("C23").Value = 1
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 8
$ws.Range("H23").Value = -12.5
$ws.Range("I23").Value = 24
$ws.Range("J23").Value = 32
$ws.Range("K23").Value = -25
$ws.Range("L23").Value = -36.842105263157
$ws.Range("M23").Value = -20

# Row 24
$ws.Range("C24").Value = 22
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 15.789473684210
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 109
$ws.Range("H24").Value = -17.431192660550
$ws.Range("I24").Value = 384
$ws.Range("J24").Value = 539
$ws.Range("K24").Value = -28.756957328385
$ws.Range("L24").Value = -33.678756476683
$ws.Range("M24").Value = 31.506849315068

# Row 25
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 71.428571428571
$ws.Range("F25").Value = 54
$ws.Range("G25").Value = 71
$ws.Range("H25").Value = -23.943661971831
$ws.Range("I25").Value = 249
$ws.Range("J25").Value = 361
$ws.Range("K25").Value = -31.024930747922
$ws.Range("L25").Value = -41.822429906542

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 7
$ws.Range("E26").Value = -14.285714285714
$ws.Range("F26").Value = 19
$ws.Range("G26").Value = 24
$ws.Range("H26").Value = -20.833333333333
$ws.Range("I26").Value = 80
$ws.Range("J26").Value = 83
$ws.Range("K26").Value = -3.614457831325
$ws.Range("L26").Value = -5.882352941176
$ws.Range("M26").Value = -9.090909090909

# Row 27
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("L27").Value = 50

# Row 28
$ws.Range("J14").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("I28").Value = 14
$ws.Range("K28").Value = -17.647058823529
$ws.Range("L28").Value = 7.692307692307

# Row 31
$ws.Range("D31").Value = 2
$ws.Range("F31").Value = 3
$ws.Range("G31").Value = 3
$ws.Range("H31").Value = 0
$ws.Range("I31").Value = 9
$ws.Range("J31").Value = 7
$ws.Range("K31").Value = 28.571428571428
$ws.Range("L31").Value = 50

$excel.CutCopyMode = 0